# Append 8 new GSW box-score rows (rows 136-143) to Sheet1, mirroring the
# layout/style of the existing data rows (row 2 .. row 135).
#
# Column layout:
#  A: game index (numeric, bold/centered/bordered "index" style copied from
#     the existing column-A cells)
#  B: TEAM, C: OPP, D: STATUS, E: DATE, F: MIN  (text)
#  G..Y: box-score numeric stats

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use the formatting of the last existing data row (135) as the template for
# the new rows' styling (column A keeps its bold/bordered "index" look, the
# rest stay on the default style) so the appended rows look like the ones
# the data-generation script already produced.
$ws.Range("A135:Y135").Copy()
$ws.Range("A136:Y143").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$newRows = @(
    @(134, "DEN", "GSW", "away", "2025-03-17", "240:00", 43, 92, 0.467, 14, 39, 0.359, 14, 18, 0.778, 19, 34, 53, 32, 12, 7, 20, 24, 114, 9),
    @(135, "GSW", "DEN", "home", "2025-03-17", "240:00", 41, 84, 0.488, 8, 33, 0.242, 15, 27, 0.556, 11, 27, 38, 24, 16, 10, 20, 21, 105, -9),
    @(136, "MIL", "GSW", "away", "2025-03-18", "240:00", 31, 79, 0.392, 12, 38, 0.316, 19, 27, 0.704, 4, 30, 34, 19, 9, 4, 9, 15, 93, -11),
    @(137, "GSW", "MIL", "home", "2025-03-18", "240:00", 32, 77, 0.416, 17, 44, 0.386, 23, 26, 0.885, 9, 43, 52, 25, 7, 6, 18, 24, 104, 11),
    @(138, "TOR", "GSW", "away", "2025-03-20", "240:00", 44, 78, 0.5639999999999999, 11, 26, 0.423, 15, 19, 0.789, 7, 37, 44, 32, 7, 4, 23, 24, 114, -3),
    @(139, "GSW", "TOR", "home", "2025-03-20", "240:00", 38, 90, 0.422, 19, 47, 0.404, 22, 28, 0.786, 11, 27, 38, 30, 13, 5, 14, 15, 117, 3),
    @(140, "GSW", "ATL", "away", "2025-03-22", "240:00", 39, 84, 0.464, 16, 38, 0.421, 21, 26, 0.8080000000000001, 6, 32, 38, 28, 8, 4, 9, 16, 115, -9),
    @(141, "ATL", "GSW", "home", "2025-03-22", "240:00", 49, 86, 0.57, 15, 36, 0.417, 11, 18, 0.611, 10, 36, 46, 37, 5, 3, 14, 23, 124, 9)
)

$startRow = 136
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $vals = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $vals[0]
    $ws.Cells.Item($r, 2).Value = $vals[1]
    $ws.Cells.Item($r, 3).Value = $vals[2]
    $ws.Cells.Item($r, 4).Value = $vals[3]

    # DATE column: the literal "YYYY-MM-DD" text would otherwise be
    # auto-parsed into a date serial by plain Value assignment. Write it as
    # a formula that yields the literal string, then freeze it to a plain
    # value in place - this keeps the cell a genuine text cell (matching the
    # source data's "YYYY-MM-DD" string rows) without leaving behind an
    # extra quote-prefixed / text-number-format style.
    $dc = $ws.Cells.Item($r, 5)
    $dc.Formula = '="' + $vals[4] + '"'
    $dc.Copy()
    $dc.PasteSpecial(-4163)
    $excel.CutCopyMode = 0

    $ws.Cells.Item($r, 6).Value = $vals[5]

    for ($c = 7; $c -le 25; $c++) {
        $ws.Cells.Item($r, $c).Value = $vals[$c - 1]
    }
}

$ws.Range("A1:Y143").Select()
$ws.Range("A1").Select()
